$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.541.79'
$ws.Range("E2").Value = '  +0.26%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.828.68'
$ws.Range("E3").Value = '  +0.14%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.11'
$ws.Range("E5").Value = '  +0.55%  '

$ws.Range("E6").Value = '  +0.03%  '

$ws.Range("E7").Value = '  +0.67%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3896'
$ws.Range("E8").Value = '  -1.11%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08421'
$ws.Range("E9").Value = '  +9.35%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.124'
$ws.Range("E10").Value = '  +1.09%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.94'
$ws.Range("E11").Value = '  +0.01%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.436'
$ws.Range("E12").Value = '  +2.51%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.26'
$ws.Range("E13").Value = '  +1.29%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.003'
$ws.Range("E14").Value = '  -0.03%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.538'
$ws.Range("E15").Value = '  -0.35%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.824.96'
$ws.Range("E16").Value = '  -0.04%  '

$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '94.10'
$ws.Range("E17").Value = '  +0.51%  '

$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001130'
$ws.Range("E18").Value = '  +4.66%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06627'
$ws.Range("E19").Value = '  -0.09%  '

$ws.Range("E20").Value = '  +0.63%  '

$ws.Range("E21").Value = '  +0.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.088'
$ws.Range("E22").Value = '  +0.54%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.574.83'
$ws.Range("E23").Value = '  +0.28%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.44'
$ws.Range("E24").Value = '  +2.69%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.278'
$ws.Range("E25").Value = '  +1.62%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '21.17'
$ws.Range("E26").Value = '  +2.83%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '159.59'
$ws.Range("E27").Value = '  +1.63%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.034.38'
$ws.Range("E28").Value = '  +0.04%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.427'
$ws.Range("E29").Value = '  -0.25%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.80'
$ws.Range("E30").Value = '  +0.66%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1101'
$ws.Range("E31").Value = '  +0.21%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.099'
$ws.Range("E32").Value = '  -2.55%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.07622'
$ws.Range("E33").Value = '  +6.35%  '

$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.747'
$ws.Range("E34").Value = '  +1.67%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.669'
$ws.Range("E35").Value = '  +0.51%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2227'
$ws.Range("E36").Value = '  -0.17%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02378'
$ws.Range("E37").Value = '  +2.26%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.270'
$ws.Range("E38").Value = '  +2.07%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.796'
$ws.Range("E39").Value = '  -2.12%  '

$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6388'
$ws.Range("E40").Value = '  +2.45%  '

$ws.Range("B41").Value = 'Aptos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.47'
$ws.Range("E41").Value = '  +1.90%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.192'
$ws.Range("E42").Value = '  +0.33%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.401'
$ws.Range("E43").Value = '  +0.53%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.66'
$ws.Range("E44").Value = '  +1.53%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6048'
$ws.Range("E45").Value = '  +2.65%  '

$ws.Range("E46").Value = '  +2.01%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '127.14'
$ws.Range("E47").Value = '  +2.20%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.997'
$ws.Range("E48").Value = '  +1.06%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.206'
$ws.Range("E49").Value = '  +2.02%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06987'
$ws.Range("E50").Value = '  +0.66%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '74.85'
$ws.Range("E51").Value = '  +1.47%  '
